# Updated symbol list on Thu Dec 15 03:42:50 UTC 2022 with GitHub Actions
#
# The "Price" column (D) stores numeric-looking quotes as plain text
# (that's how the source sheet is generated), so a leading apostrophe is
# used to force each assignment to stay text instead of being coerced to
# a number by Excel.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value  = "'265.85"
$ws.Range("D3").Value  = "'22.57"
$ws.Range("D4").Value  = "'6.278"
$ws.Range("D5").Value  = "'0.06160"
$ws.Range("D6").Value  = "'3.572"
$ws.Range("D8").Value  = "'1.348"
$ws.Range("D9").Value  = "'0.8278"
$ws.Range("D11").Value = "'0.1589"
$ws.Range("D12").Value = "'0.08189"
$ws.Range("D13").Value = "'0.03424"
$ws.Range("D14").Value = "'0.03178"
$ws.Range("D15").Value = "'0.09248"
$ws.Range("D16").Value = "'3.887"
$ws.Range("D17").Value = "'0.001702"
$ws.Range("D18").Value = "'0.04883"
$ws.Range("D19").Value = "'0.006220"
$ws.Range("D20").Value = "'0.005274"
$ws.Range("D24").Value = "'2.321"
$ws.Range("D26").Value = "'0.1237"
$ws.Range("D27").Value = "'0.0002679"
$ws.Range("D40").Value = "'0.04632"
$ws.Range("D41").Value = "'0.006953"
$ws.Range("D42").Value = "'0.1138"
$ws.Range("D43").Value = "'0.003130"
$ws.Range("E44").Value = "43LocalTradersLCT"
$ws.Range("D45").Value = "'0.00006137"
$ws.Range("D47").Value = "'0.6998"
$ws.Range("E47").Value = "46CoinbaseStockTokenCOINWorstin24h"
$ws.Range("D48").Value = "'0.1931"
$ws.Range("D49").Value = "'0.00002099"
